# Auto-generated edit script: updates the cryptos price/volume table
# to the refreshed values captured by the GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.561.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5262'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3296'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06752'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.45'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7779'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07649'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.825.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.057'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007911'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.581.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.076.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.615'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.736'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.010'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.372'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.648'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.237'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.196'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08790'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04858'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.145'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.857'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7101'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.109'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01815'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.222'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.34%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4970'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '114.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9041'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.076'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.826'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.0000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4303'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1294'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.180'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05926'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.43'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.440'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.91%  '
